# Second Commit of the day
# Applies the HackathonData.xlsx edits described by the diff:
#  - Several placeholder "Test data N" numeric columns get reset to "NA"
#  - TC_003's test case name gains a longer suffix
#  - TC_010's "Test data 3"/"Test data 4" columns get populated with real URLs
#  - Two brand new test case rows (TC_013 / TC_014) are appended
#  - Row 11 grows taller to accommodate the newly-populated URL text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: some columns in this sheet use a "quote-prefixed" text style (the
# cell already forces numbers such as "2" or "27" to be stored as text). When
# COM sets a plain numeric-looking string, Excel converts it to a real number
# and drops that style. Prefixing the value with a leading apostrophe keeps
# it text and preserves the existing quote-prefix style, matching how this
# sheet already stores its numbers-as-text.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$value)
    $range.Value = "'" + $value
}

# --- Row 2 (TC_001) -------------------------------------------------------
Set-TextValue $ws.Range("E2") "NA"
Set-TextValue $ws.Range("F2") "NA"
Set-TextValue $ws.Range("G2") "NA"
Set-TextValue $ws.Range("H2") "NA"
Set-TextValue $ws.Range("I2") "NA"
Set-TextValue $ws.Range("J2") "NA"

# --- Row 3 (TC_002) --------------------------------------------------------
Set-TextValue $ws.Range("J3") "NA"

# --- Row 4 (TC_003) --------------------------------------------------------
$ws.Range("C4").Value = "verifyIfCourseCardsArePresentAfterApplyingFilters"
Set-TextValue $ws.Range("J4") "NA"

# --- Row 9 (TC_008) --------------------------------------------------------
Set-TextValue $ws.Range("G9") "NA"
Set-TextValue $ws.Range("H9") "NA"
Set-TextValue $ws.Range("I9") "NA"
Set-TextValue $ws.Range("J9") "NA"

# --- Row 10 (TC_009) --------------------------------------------------------
Set-TextValue $ws.Range("F10") "NA"
Set-TextValue $ws.Range("G10") "NA"
Set-TextValue $ws.Range("H10") "NA"
Set-TextValue $ws.Range("I10") "NA"
Set-TextValue $ws.Range("J10") "NA"

# --- Row 11 (TC_010) - now gets the real page title / URL data -----------
Set-TextValue $ws.Range("F11") "Language Learning Online Courses | Coursera"
Set-TextValue $ws.Range("G11") "Coursera Online Course Catalog by Topic and Skill | Coursera"
Set-TextValue $ws.Range("H11") "NA"
Set-TextValue $ws.Range("I11") "NA"
Set-TextValue $ws.Range("J11") "NA"
$ws.Rows.Item(11).RowHeight = 101.5

# --- Row 12 (TC_011) --------------------------------------------------------
Set-TextValue $ws.Range("G12") "NA"
Set-TextValue $ws.Range("H12") "NA"
Set-TextValue $ws.Range("I12") "NA"
Set-TextValue $ws.Range("J12") "NA"

# --- Row 13 (TC_012) --------------------------------------------------------
Set-TextValue $ws.Range("H13") "NA"
Set-TextValue $ws.Range("I13") "NA"
Set-TextValue $ws.Range("J13") "NA"

# --- Row 14 (new TC_013) ----------------------------------------------------
$ws.Range("A14").Value = "TS002"
$ws.Range("B14").Value = "TC_013"
$ws.Range("C14").Value = "getCountOfLanguages"
Set-TextValue $ws.Range("D14") "Explore Categories"
Set-TextValue $ws.Range("E14") "Language Learning"
$ws.Range("F14").Value = "Language"
Set-TextValue $ws.Range("G14") "27"
Set-TextValue $ws.Range("H14") "8"
Set-TextValue $ws.Range("I14") "NA"
Set-TextValue $ws.Range("J14") "NA"
$ws.Rows.Item(14).RowHeight = 43.5

# --- Row 15 (new TC_014) ----------------------------------------------------
$ws.Range("A15").Value = "TS002"
$ws.Range("B15").Value = "TC_014"
$ws.Range("C15").Value = "getCountOfLevels"
Set-TextValue $ws.Range("D15") "Explore Categories"
Set-TextValue $ws.Range("E15") "Language Learning"
$ws.Range("F15").Value = "Level"
Set-TextValue $ws.Range("G15") "4"
Set-TextValue $ws.Range("H15") "3"
Set-TextValue $ws.Range("I15") "NA"
Set-TextValue $ws.Range("J15") "NA"
$ws.Rows.Item(15).RowHeight = 43.5

# --- Update the sheet view to match where the author last left it ---------
$ws.Range("E12").Select()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 2
